# edit.ps1 - applies the "Modify Loading Gif, Add Sentence Popup" commit
# to ECommuBook-APIs.docx via Word COM interop.
#
# Summary of edits (see accompanying diff):
#   1. Merge the "Search" + ":" runs into a single "Search:" run.
#   2. Append a red " (Need Server Side)" run after "Voice conversion".
#   3. Change "Popup edit (Need Server side)" -> "Popup edit (Complete)".
#   4. Move the "_GoBack" bookmark: it used to sit (alone) in the empty
#      paragraph just above "Keyword search and Record (Need Server side)";
#      now that paragraph holds the "Keyword search..." text instead, and
#      the bookmark moves down to sit at the start of the (relocated)
#      "Download Page: load time give GIF" paragraph.
#   5. Delete the now-redundant "Image to GIF function" paragraph and the
#      old "Download Page: load time give GIF" paragraph (its text lives
#      on in the paragraph that used to hold the bookmark).

$d = $word.ActiveDocument

# Helper: locate the Paragraph object containing document position $pos.
# (Range.Paragraphs.Item(1) does not reliably scope to the owning range in
# this host, so we resolve paragraphs via Document.Paragraphs instead.)
function Get-ParagraphAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Start -le $pos -and $pos -lt $cand.Range.End) {
            return $cand
        }
    }
    return $doc.Paragraphs.Item($doc.Paragraphs.Count)
}

# Helper: same as above but returns the 1-based index instead of the
# Paragraph object. Cached Paragraph/Range object references do not update
# their Start/End after an earlier paragraph is deleted in this host, so
# structural edits below re-resolve paragraphs by index after every delete.
function Get-ParagraphIndexAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Start -le $pos -and $pos -lt $cand.Range.End) {
            return $i
        }
    }
    return $doc.Paragraphs.Count
}

# ---------------------------------------------------------------------
# 1. "Search" + ":" -> single run "Search:"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Search:", $true, $false, $false, $false, $false, $true, 1, $false, "Search:", 2)

# ---------------------------------------------------------------------
# 2. "Voice conversion" -> "Voice conversion" + red " (Need Server Side)"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Voice conversion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" (Need Server Side)")
$rng.Font.Color = 255

# ---------------------------------------------------------------------
# 3. "Popup edit (Need Server side)" -> "Popup edit (Complete)"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Popup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$popupPara = Get-ParagraphAt $d $rng.Start
$rng2 = $popupPara.Range
$rng2.Find.Execute(" (Need Server side)", $true, $false, $false, $false, $false, $true, 1, $false, " (Complete)", 2)

# ---------------------------------------------------------------------
# 4 & 5. Bookmark relocation + paragraph cleanup around the GIF bullets
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Sentence popup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentenceIdx = Get-ParagraphIndexAt $d $rng.Start

$bookmarkIdx = $sentenceIdx + 1   # currently empty, holds the "_GoBack" bookmark
$keywordIdx  = $sentenceIdx + 2   # "Keyword search and Record (Need Server side)"
$gifIdx      = $sentenceIdx + 3   # "Image to GIF function"
# ($sentenceIdx + 4 is "Download Page: load time give GIF" - resolved by
#  index, not cached, once the deletes above it have happened)

# Remove the bookmark that currently lives in the (empty) $bookmarkIdx
# paragraph, then give that paragraph the relocated
# "Keyword search and Record (Need Server side)" text.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$bmRng = $d.Paragraphs.Item($bookmarkIdx).Range
$bmRng.End = $bmRng.End - 1
$bmRng.InsertBefore("Keyword search and Record (Need Server side)")

# Replace the $keywordIdx paragraph's old text with
# "Download Page: load time give GIF" and re-add the bookmark at its start.
$kwRng = $d.Paragraphs.Item($keywordIdx).Range
$kwRng.End = $kwRng.End - 1
$kwRng.Text = "Download Page: load time give GIF"
$kwStart = $d.Paragraphs.Item($keywordIdx).Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($kwStart, $kwStart))

# Delete the obsolete "Image to GIF function" paragraph...
$d.Paragraphs.Item($gifIdx).Range.Delete()
# ...which shifts the old "Download Page: load time give GIF" paragraph
# (whose content has already been relocated above) down into $gifIdx; drop it
# too.
$d.Paragraphs.Item($gifIdx).Range.Delete()
